$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes existing rows 32:151 down to 33:152)
$ws.Rows("32:32").Insert()

# Populate the new row 32 with the new weekly record
$ws.Cells.Item(32, 1).Value = 5
$ws.Cells.Item(32, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(32, 3).Value = "Maule"
$ws.Cells.Item(32, 4).Value = 44624
$ws.Cells.Item(32, 5).Value = 7
$ws.Cells.Item(32, 6).Value = 100112031
$ws.Cells.Item(32, 7).Value = "Poroto verde"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 150
$ws.Cells.Item(32, 11).Value = 30000
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = 30000
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 1200
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Make sure the new D32 cell carries the same date-style formatting (style index 2)
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(33, 4).NumberFormat
